$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "kvw5270"
$ws.Range("B2").Value = "03/24/2020 01:00:57"
